$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 245
$ws.Range("F3").Value = 842
$ws.Range("F4").Value = 551
$ws.Range("F5").Value = 2285
$ws.Range("F6").Value = 1361
$ws.Range("F9").Value = 1151
$ws.Range("F10").Value = 494
$ws.Range("F11").Value = 3016
$ws.Range("F14").Value = 1102
$ws.Range("F17").Value = 232
$ws.Range("F19").Value = 1083
$ws.Range("F20").Value = 1083
$ws.Range("F22").Value = 529
$ws.Range("F23").Value = 174
$ws.Range("F25").Value = 224
$ws.Range("F29").Value = 834
$ws.Range("F32").Value = 20
$ws.Range("F33").Value = 1039
$ws.Range("F34").Value = 5040
$ws.Range("F35").Value = 500
$ws.Range("F36").Value = 239
$ws.Range("F37").Value = 126

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 22
$ws.Range("F11").Value = 195
$ws.Range("F14").Value = 9
$ws.Range("F17").Value = 65
$ws.Range("F22").Value = 301
$ws.Range("F24").Value = 49
$ws.Range("F28").Value = 668
$ws.Range("F34").Value = 60
$ws.Range("F37").Value = 436
$ws.Range("F43").Value = 753

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 639
$ws.Range("F6").Value = 408

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 639
$ws.Range("F3").Value = 245
$ws.Range("F5").Value = 842
$ws.Range("F7").Value = 551
$ws.Range("F8").Value = 2285
$ws.Range("F9").Value = 1361
$ws.Range("F12").Value = 1151
$ws.Range("F14").Value = 195
$ws.Range("F15").Value = 498
$ws.Range("F16").Value = 3016
$ws.Range("F19").Value = 1102
$ws.Range("F21").Value = 408
$ws.Range("F23").Value = 232
$ws.Range("F25").Value = 1083
$ws.Range("F26").Value = 1083
$ws.Range("F29").Value = 529
$ws.Range("F30").Value = 301
$ws.Range("F31").Value = 174
$ws.Range("F32").Value = 224
$ws.Range("F33").Value = 49
$ws.Range("F37").Value = 668
$ws.Range("F38").Value = 834
$ws.Range("F43").Value = 1039
$ws.Range("F44").Value = 5040
$ws.Range("F45").Value = 60
$ws.Range("F46").Value = 500
$ws.Range("F47").Value = 436
$ws.Range("F48").Value = 436
$ws.Range("F49").Value = 239
$ws.Range("F51").Value = 753
